# ---------------------------------------------------------------------------
# Commit: "complate!!-> scrapping whole ipl"
#
# Rebuilds "M Shahrukh Khan.xlsx": renames the default "Sheet1" tab to the
# batter's name, adds a leading "matchNo" column, and expands the single
# sample row into the full set of scraped innings (rows 2-11), A1:M11.
# ---------------------------------------------------------------------------

# Full target grid: row 1 is the header, rows 2-11 are one innings each.
# Every value is written as text (scores like "8" or "80.00" must stay
# strings, matching the source scraper's output) - an empty string marks a
# "not out" innings with no recorded dismissal ("states").
$data = @(
    @('matchNo','teamName','batterName','states','runs','balls','fours','sixes','sr','opponentTeamName','venue','date','result'),
    @('53rd','Punjab Kings','M Shahrukh Khan','c Bravo b Chahar','8','10','0','1','80.00','Chennai Super Kings','Dubai (DSC)','October 07','Punjab Kings won by 6 wickets (with 42 balls remaining)'),
    @('45th','Punjab Kings','M Shahrukh Khan','','22','9','1','2','244.44','Kolkata Knight Riders','Dubai (DSC)','October 01','Punjab Kings won by 5 wickets (with 3 balls remaining)'),
    @('8th','Punjab Kings','M Shahrukh Khan','c Jadeja b Curran','47','36','4','2','130.55','Chennai Super Kings','Wankhede','April 16','Super Kings won by 6 wickets (with 26 balls remaining)'),
    @('11th','Punjab Kings','M Shahrukh Khan','','15','5','2','1','300.00','Delhi Capitals','Wankhede','April 18','Capitals won by 6 wickets (with 10 balls remaining)'),
    @('48th','Punjab Kings','M Shahrukh Khan','run out (Patel)','16','11','1','1','145.45','Royal Challengers Bangalore','Sharjah','October 03','RCB won by 6 runs'),
    @('21st','Punjab Kings','M Shahrukh Khan','c Morgan b Prasidh Krishna','13','14','0','1','92.85','Kolkata Knight Riders','Ahmedabad','April 26','KKR won by 5 wickets (with 20 balls remaining)'),
    @('26th','Punjab Kings','M Shahrukh Khan','b Chahal','0','3','0','0','0.00','Royal Challengers Bangalore','Ahmedabad','April 30','Punjab Kings won by 34 runs'),
    @('29th','Punjab Kings','M Shahrukh Khan','c Hetmyer b Avesh Khan','4','5','0','0','80.00','Delhi Capitals','Ahmedabad','May 02','Capitals won by 7 wickets (with 14 balls remaining)'),
    @('14th','Punjab Kings','M Shahrukh Khan','c Abhishek Sharma b Ahmed','22','17','0','2','129.41','Sunrisers Hyderabad','Chennai','April 21','Sunrisers won by 9 wickets (with 8 balls remaining)'),
    @('4th','Punjab Kings','M Shahrukh Khan','','6','4','1','0','150.00','Rajasthan Royals','Wankhede','April 12','Punjab Kings won by 4 runs')
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sheet1" -> "M Shahrukh Khan"
$ws.Name = "M Shahrukh Khan"

$numRows = $data.Count
$numCols = $data[0].Count

# Wipe whatever the sheet held (the old 2x12 sample table) before laying the
# new 11x13 table down, so nothing from the previous shape lingers.
$ws.Cells.Clear()

# Every cell in this sheet is scraped text, never a "real" number - force
# Text format over the whole target range first so values such as "0",
# "8" or "80.00" are kept exactly as typed instead of being coerced to
# numbers.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($numRows, $numCols))
$fullRange.NumberFormat = "@"

for ($r = 0; $r -lt $numRows; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
